{"js": "// Support for common specializations decisions\n//\n// 1) \"\u0398\u03ad\u03bc\u03b1: \u00ab\u0388\u03b3\u03ba\u03c1\u03b9\u03c3\u03b7 \u03b1\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7\u03c2 \u03b4\u03b9\u03ac\u03b8\u03b5\u03c3\u03b7\u03c2 \u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03ce\u03bd \u03c4\u03b7\u03c2 ${local_directorate}\u00bb\"\n//    -> the first ${local_directorate} placeholder (the heading/subject line)\n//       becomes ${local_directorate_genitive}\n// 2) \"... \u03b3\u03b9\u03b1 \u03c4\u03bf\u03c5\u03c2 \u03ba\u03ac\u03c4\u03c9\u03b8\u03b9 \u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03bf\u03cd\u03c2 , \u03c9\u03c2 \u03b5\u03be\u03ae\u03c2:\"\n//    -> drop the stray space before the comma:\n//       \"... \u03b3\u03b9\u03b1 \u03c4\u03bf\u03c5\u03c2 \u03ba\u03ac\u03c4\u03c9\u03b8\u03b9 \u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03bf\u03cd\u03c2, \u03c9\u03c2 \u03b5\u03be\u03ae\u03c2:\"\n\nconst body = context.document.body;\n\n// --- Change 1: ${local_directorate} -> ${local_directorate_genitive} -----\n// Only the placeholder immediately followed by \"\u00bb\" (the subject/heading use)\n// is renamed; the later \"${local_directorate}\" before \"\u039a\u039f\u0399\u039d\u039f\u03a0\u039f\u0399\u0397\u03a3\u0397\" (used for\n// the CC/notification block) is left untouched, matching the source edit.\nconst headingResults = body.search(\"${local_directorate}\u00bb\", { matchCase: true });\nheadingResults.load(\"text\");\nawait context.sync();\n\nif (headingResults.items.length === 0) {\n  throw new Error(\"Could not find the '${local_directorate}\u00bb' placeholder to update.\");\n}\n\nheadingResults.items[0].insertText(\"${local_directorate_genitive}\u00bb\", \"Replace\");\nawait context.sync();\n\n// --- Change 2: remove the extra space before the comma -------------------\nconst spaceCommaResults = body.search(\"\u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03bf\u03cd\u03c2 , \u03c9\u03c2 \u03b5\u03be\u03ae\u03c2:\", { matchCase: true });\nspaceCommaResults.load(\"text\");\nawait context.sync();\n\nif (spaceCommaResults.items.length === 0) {\n  throw new Error(\"Could not find the '\u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03bf\u03cd\u03c2 , \u03c9\u03c2 \u03b5\u03be\u03ae\u03c2:' text to update.\");\n}\n\nspaceCommaResults.items[0].insertText(\"\u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03bf\u03cd\u03c2, \u03c9\u03c2 \u03b5\u03be\u03ae\u03c2:\", \"Replace\");\nawait context.sync();\n", "ps1": "# Support for common specializations decisions\n#\n# 1) \"\u0398\u03ad\u03bc\u03b1: \u00ab\u0388\u03b3\u03ba\u03c1\u03b9\u03c3\u03b7 \u03b1\u03c0\u03cc\u03c6\u03b1\u03c3\u03b7\u03c2 \u03b4\u03b9\u03ac\u03b8\u03b5\u03c3\u03b7\u03c2 \u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03ce\u03bd \u03c4\u03b7\u03c2 ${local_directorate}\u00bb\"\n#    -> rename the heading/subject-line placeholder to\n#       ${local_directorate_genitive}\n# 2) \"... \u03b3\u03b9\u03b1 \u03c4\u03bf\u03c5\u03c2 \u03ba\u03ac\u03c4\u03c9\u03b8\u03b9 \u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03bf\u03cd\u03c2 , \u03c9\u03c2 \u03b5\u03be\u03ae\u03c2:\"\n#    -> drop the stray space before the comma.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: ${local_directorate} -> ${local_directorate_genitive} -----\n# Only the placeholder immediately followed by \"\u00bb\" (the subject/heading use)\n# is renamed; the other \"${local_directorate}\" occurrence later in the\n# document (before \"\u039a\u039f\u0399\u039d\u039f\u03a0\u039f\u0399\u0397\u03a3\u0397\") is left untouched, matching the source edit.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = '${local_directorate}\u00bb'\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = '${local_directorate_genitive}\u00bb'\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute(\n    $find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards,\n    $false, $false, $find.Forward, $find.Wrap, $false,\n    $find.Replacement.Text, 1  # wdReplaceOne\n) | Out-Null\n\n# --- Change 2: remove the extra space before the comma -------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = '\u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03bf\u03cd\u03c2 , \u03c9\u03c2 \u03b5\u03be\u03ae\u03c2:'\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = '\u03b5\u03ba\u03c0\u03b1\u03b9\u03b4\u03b5\u03c5\u03c4\u03b9\u03ba\u03bf\u03cd\u03c2, \u03c9\u03c2 \u03b5\u03be\u03ae\u03c2:'\n$find2.Forward = $true\n$find2.Wrap = 0  # wdFindStop\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.MatchWildcards = $false\n$find2.Execute(\n    $find2.Text, $find2.MatchCase, $find2.MatchWholeWord, $find2.MatchWildcards,\n    $false, $false, $find2.Forward, $find2.Wrap, $false,\n    $find2.Replacement.Text, 1  # wdReplaceOne\n) | Out-Null\n"}
